$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Data changes -----------------------------------------------------
# Add the new recursion note (Sheet2!D3), then the newly-worked-out
# solution set that replaces the placeholder "9" at Sheet2!B10.
# (Order matters for the shared-string table append order.)
$ws2.Range("D3").Value = "297?"
$ws2.Range("B10").Value = "22,98,104,111,226"

# Sheet1!B10 (=Sheet2!B10) now holds a real solved answer instead of a
# bare recursion placeholder, so give it the same "solved" highlight
# fill the other answered cells (e.g. B9) already use.
$ws1.Range("B9").Copy()
$ws1.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths ------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 19.33
$ws1.Columns.Item(4).ColumnWidth = 18.5
$ws2.Columns.Item(2).ColumnWidth = 22.67

# --- Window / selection state --------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 23
$win.Width = 1440
$win.Height = 877

$ws2.Activate() | Out-Null
$ws2.Range("B11").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D9").Select() | Out-Null
